$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Orientation" header in D1, matching the bold header style used by A1:C1
$ws.Range("D1").Value = "Orientation"
$ws.Range("D1").Font.Bold = $true

# Four antenna/receiver rows with lat/long + orientation degrees
$ws.Range("A2").Value = "receiver_1"
$ws.Range("B2").Value = 1.234
$ws.Range("C2").Value = 5.678
$ws.Range("D2").Value = 0

$ws.Range("A3").Value = "receiver_2"
$ws.Range("B3").Value = 1.234
$ws.Range("C3").Value = 5.678
$ws.Range("D3").Value = 90

$ws.Range("A4").Value = "receiver_3"
$ws.Range("B4").Value = 1.234
$ws.Range("C4").Value = 5.678
$ws.Range("D4").Value = 180

$ws.Range("A5").Value = "receiver_4"
$ws.Range("B5").Value = 1.234
$ws.Range("C5").Value = 5.678
$ws.Range("D5").Value = 270

# Widen column A (Name) and D (Orientation) to fit their new content
$ws.Columns.Item(1).ColumnWidth = 13.333333333333334
$ws.Columns.Item(4).ColumnWidth = 10.5

# Leave the selection where Excel would land after typing the last value
[void]$ws.Range("D6").Select()
